$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give every staff record a unique id (column J) -- new values for rows 2-4
$ws.Range("J2").Value = "VEC-008-01-191"
$ws.Range("J3").Value = "VEC-008-04-174"
$ws.Range("J4").Value = "VEC-008-04-172"

# Column I grew to fit the longer unique-id strings, and the new column J
# needs an explicit width too (best achievable values given the engine's
# internal column-width quantization).
$ws.Columns.Item(9).ColumnWidth = 36
$ws.Columns.Item(10).ColumnWidth = 17.666666666666668

# Match the author's final selection/viewport
$ws.Range("J5").Select()
